$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 8
$ws.Range("B6").Value = "Terminal La Palmera de La Serena"
$ws.Range("C6").Value = "Coquimbo"
$ws.Range("D6").Value = 44960
$ws.Range("D6").NumberFormat = $ws.Range("D5").NumberFormat
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 100114002
$ws.Range("G6").Value = "Camote"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 400
$ws.Range("K6").Value = 19500
$ws.Range("L6").Value = 20000
$ws.Range("M6").Value = 19750
$ws.Range("N6").Value = "`$/malla 18 kilos"
$ws.Range("O6").Value = "Perú"
$ws.Range("P6").Value = 1097
$ws.Range("Q6").Value = 18
$ws.Range("R6").Value = "Hortaliza"
